$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2, 4, 6 (K unchanged = 1)
$valsOdd = @{
    'C' = 0.0001
    'E' = 1000
    'J' = 0.002
    'K' = 1
    'L' = 0.9980000257492065
    'M' = 0.9937000274658203
    'N' = 105.559
    'O' = 0.0095
    'P' = 0.0094
    'Q' = 14
    'R' = 7.54
    'S' = 0.9962999820709229
    'T' = 0.9950000047683716
    'U' = 0.9957000017166138
}

# New values for rows 3, 5 (K becomes 2)
$valsEven = @{
    'C' = 0.0001
    'E' = 1000
    'J' = 0.002
    'K' = 2
    'L' = 0.9993000030517578
    'M' = 0.9973999857902527
    'N' = 166.16
    'O' = 0.006
    'P' = 0.006
    'Q' = 21
    'R' = 7.9124
    'S' = 0.998199999332428
    'T' = 0.9947999715805054
    'U' = 0.9980999827384949
}

$rowsWithOddPattern = @(2, 4, 6)
$rowsWithEvenPattern = @(3, 5)

foreach ($r in $rowsWithOddPattern) {
    foreach ($col in $valsOdd.Keys) {
        $ws.Range("$col$r").Value = $valsOdd[$col]
    }
}

foreach ($r in $rowsWithEvenPattern) {
    foreach ($col in $valsEven.Keys) {
        $ws.Range("$col$r").Value = $valsEven[$col]
    }
}
